$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист")

$ws.Range("A52").Value = "Павликов Илья Сергеевич "
$ws.Range("A53").Value = "Верле Каролина Валерьевна (Обучение 2)"
$ws.Range("A54").Value = "Довыдович Алиса Станиславовна "
$ws.Range("A55").Value = " Шептунова Софья Денисовна"
